$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- "Personas" row (row 15 in 1-based Word table indexing) ---
# Remove the "Personas" run text, leaving the (now empty) paragraph intact.
[void]$d.Content.Find.Execute("Personas", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# The second cell of that row has 5 empty paragraphs; keep the first, drop the other four.
$cell = $t.Cell(15,3)
for ($i = $cell.Range.Paragraphs.Count; $i -ge 2; $i--) {
    $cell.Range.Paragraphs.Item($i).Range.Delete()
}

# --- "User storyboards" row (row 16) ---
# Remove the "User storyboards" run text.
[void]$d.Content.Find.Execute("User storyboards", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# That cell also had a second, now-redundant empty paragraph following the text paragraph; drop it.
$cell = $t.Cell(16,1)
$cell.Range.Paragraphs.Item($cell.Range.Paragraphs.Count).Range.Delete()

# The second cell of that row has 5 empty paragraphs; keep the last, drop the first four.
$cell = $t.Cell(16,3)
for ($i = 1; $i -le 4; $i++) {
    $cell.Range.Paragraphs.Item(1).Range.Delete()
}

# --- "Further comments:" paragraph: drop the stale lastRenderedPageBreak marker ---
$cell = $t.Cell(18,1)
$p = $cell.Range.Paragraphs.Item(1)
$p.Range.Text = "Further comments:"

# --- Tutor row: remove the "RE/" + "DW" runs, leaving an empty paragraph ---
[void]$d.Content.Find.Execute("RE/DW", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
